# Add missing input variable IOP_INTERPOLATOR (IOP_INTERPOLATION) to the
# STR_STAB and Z_JACKET sheets of template_conductor_1_operation.xlsx, and
# refresh the STR_MIX sheet's existing IOP_INTERPOLATION note to the more
# specific "operating current" wording (previously it reused the magnetic
# field description).

$wb = $excel.ActiveWorkbook

$noteText = "Flag to specify the kind of interpolator for the operating current data, if IOP_MODE = -1. Possible values: linear = liear interpolation; cubic = use spline function of third order."

# ---------------------------------------------------------------------
# STR_MIX: row 5 (IOP_INTERPOLATION) already exists - just point its
# Note/comments cell (column D) at the new, more specific description.
# ---------------------------------------------------------------------
$wsMix = $wb.Worksheets.Item("STR_MIX")
$wsMix.Cells.Item(5, 4).Value2 = $noteText
$wsMix.Activate()
$wsMix.Range("D5").Select()

# ---------------------------------------------------------------------
# STR_STAB: insert a new row 5 carrying the IOP_INTERPOLATION variable.
# ---------------------------------------------------------------------
$wsStab = $wb.Worksheets.Item("STR_STAB")
$wsStab.Rows.Item(5).Insert()
$wsStab.Rows.Item(5).RowHeight = 45
$wsStab.Cells.Item(5, 1).Value2 = "IOP_INTERPOLATION"
$wsStab.Cells.Item(5, 2).Value2 = "-"
$wsStab.Cells.Item(5, 3).Value2 = "string"
$wsStab.Cells.Item(5, 4).Value2 = $noteText
$wsStab.Cells.Item(5, 5).Value2 = "linear"
$wsStab.Cells.Item(5, 1).Style = "Normal"
$wsStab.Activate()
$wsStab.Range("D5").Select()

# ---------------------------------------------------------------------
# Z_JACKET: same new row, but this sheet ends up the active tab and the
# selection on it lands on K14 (per the authored workbook).
# ---------------------------------------------------------------------
$wsJacket = $wb.Worksheets.Item("Z_JACKET")
$wsJacket.Rows.Item(5).Insert()
$wsJacket.Rows.Item(5).RowHeight = 45
$wsJacket.Cells.Item(5, 1).Value2 = "IOP_INTERPOLATION"
$wsJacket.Cells.Item(5, 2).Value2 = "-"
$wsJacket.Cells.Item(5, 3).Value2 = "string"
$wsJacket.Cells.Item(5, 4).Value2 = $noteText
$wsJacket.Cells.Item(5, 5).Value2 = "linear"
$wsJacket.Activate()
$wsJacket.Range("K14").Select()
